# Rename text (translate Vietnamese recommendation text to English) and
# update the timestamp identifiers in column H.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = 45786.07619457495

$rows = @(
    @{ Row = 2;  Safety = "1552 → 1241"; Reorder = "189523 → 170570" },
    @{ Row = 3;  Safety = "501 → 400";   Reorder = "40167 → 36150" },
    @{ Row = 4;  Safety = "550 → 440";   Reorder = "8075 → 7267" },
    @{ Row = 5;  Safety = "146 → 116";   Reorder = "7536 → 6782" },
    @{ Row = 6;  Safety = "152 → 121";   Reorder = "6341 → 5706" },
    @{ Row = 7;  Safety = "70 → 56";     Reorder = "5684 → 5115" },
    @{ Row = 8;  Safety = "894 → 715";   Reorder = "3548 → 3193" },
    @{ Row = 9;  Safety = "0 → 0";       Reorder = "3749 → 3374" },
    @{ Row = 10; Safety = "283 → 226";   Reorder = "2317 → 2085" },
    @{ Row = 11; Safety = "0 → 0";       Reorder = "2035 → 1831" }
)

foreach ($item in $rows) {
    $r = $item.Row
    $text = "Reduce Safety Stock from $($item.Safety) and Reorder Point from $($item.Reorder) to save costs."
    $ws.Range("B$r").Value = $text
    $ws.Range("H$r").Value = $newTimestamp
}
